# DaySale report update:
#  - Insert 6 new product rows (table grows from 4 to 10 rows)
#  - Update existing rows' values (name/balance/price/saleprice/transactions)
#  - Recompute the total and move the total/footer rows down
#  - Update the generated timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: relocate the trailing "total" row and the footer row out of the
# way (they move from rows 11/12 down to rows 17/18 once 6 product rows are
# inserted above them).
# ---------------------------------------------------------------------------

# Footer row: 12 -> 18 (values then formats, so existing styles are reused)
$ws.Range("A12:Q12").Copy()
$ws.Range("A18:Q18").PasteSpecial(-4163)
$ws.Range("A12:Q12").Copy()
$ws.Range("A18:Q18").PasteSpecial(-4122)
$ws.Rows.Item(18).RowHeight = $ws.Rows.Item(12).RowHeight
$ws.Range("A12:Q12").ClearContents()

# Total row: only P11:Q11 hold data -> P17:Q17
$ws.Range("P11:Q11").Copy()
$ws.Range("P17:Q17").PasteSpecial(-4163)
$ws.Range("P11:Q11").Copy()
$ws.Range("P17:Q17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = $ws.Rows.Item(11).RowHeight

# Fix up the merges: drop the old ones, add the new ones
$ws.Range("A12:F12").UnMerge()
$ws.Range("G12:I12").UnMerge()
$ws.Range("K12:Q12").UnMerge()
$ws.Range("P11:Q11").UnMerge()

$ws.Range("A18:F18").Merge()
$ws.Range("G18:I18").Merge()
$ws.Range("K18:Q18").Merge()
$ws.Range("P17:Q17").Merge()

# Update the total to match the new set of 10 rows
$ws.Range("P17").Value = "597.94"

# Update the generated timestamp text in the footer
$ws.Range("A18").Value = "Tuesday, 15 July, 2025 9:42 AM"

# ---------------------------------------------------------------------------
# Step 2: build the 6 new product rows (11-16) using row 10 as the
# formatting template, then fill every row 7-16 with the final data set.
# ---------------------------------------------------------------------------

$newRows = @(11, 12, 13, 14, 15, 16)
foreach ($r in $newRows) {
    $ws.Range("A11:Q11").ClearContents()
    $ws.Range("A10:Q10").Copy()
    $ws.Range("A$r`:Q$r").PasteSpecial(-4163)
    $ws.Range("A10:Q10").Copy()
    $ws.Range("A$r`:Q$r").PasteSpecial(-4122)
    $ws.Range("A$r`:B$r").Merge()
    $ws.Range("C$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
    $ws.Range("N$r`:O$r").Merge()
}

# Row heights, matching the original report's auto-fit values
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75
$ws.Rows.Item(14).RowHeight = 25.5
$ws.Rows.Item(15).RowHeight = 24.75
$ws.Rows.Item(16).RowHeight = 25.5

# ---------------------------------------------------------------------------
# Step 3: populate the table - 10 product rows, sorted by product name.
# ---------------------------------------------------------------------------

$data = @(
    @{ Row=7;  Idx=1;  Name="AMARYL 4 MG 30 TABS";               Bal="0:1"; Price="108.00"; Sale="108.0000"; Trans="1:0" },
    @{ Row=8;  Idx=2;  Name="BABY RELIEF 12.5MG 5 SUPP.";        Bal="1:0"; Price="36.00";  Sale="36.0000";  Trans="1:0" },
    @{ Row=9;  Idx=3;  Name="CONGESTAL 20 TABS";                 Bal="3:1"; Price="50.00";  Sale="25.0000";  Trans="0:1" },
    @{ Row=10; Idx=4;  Name="DEXAZONE 0.5MG 60 TAB";             Bal="1:0"; Price="51.00";  Sale="16.8300";  Trans="0:1" },
    @{ Row=11; Idx=5;  Name="ERASTAPEX CO 5/20MG 30 F.C. TABS";  Bal="1:0"; Price="114.00"; Sale="114.0000"; Trans="1:0" },
    @{ Row=12; Idx=6;  Name="FAROVIGA 100MG 12 F.C.TAB.";        Bal="2:2"; Price="108.00"; Sale="35.6400";  Trans="0:4" },
    @{ Row=13; Idx=7;  Name="OMEGA-3 PLUS 30 CAPS.";             Bal="2:2"; Price="135.00"; Sale="44.5500";  Trans="0:1" },
    @{ Row=14; Idx=8;  Name="SELENIUM-ACE 30 TABS";              Bal="1:1"; Price="130.00"; Sale="42.9000";  Trans="0:1" },
    @{ Row=15; Idx=9;  Name="SPASMO-DIGESTIN 30 TABS.";          Bal="4:0"; Price="78.00";  Sale="78.0000";  Trans="1:0" },
    @{ Row=16; Idx=10; Name="ZYROVAZET 10/20MG 30 F.C. TABLETS"; Bal="1:0"; Price="294.00"; Sale="97.0200";  Trans="0:1" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Idx           # A - index
    $ws.Range("C$r").Value = $item.Name               # C - product name
    $ws.Range("H$r").Value = $item.Bal                # H - current balance
    $ws.Range("L$r").Value = "1"                       # L - reorder limit
    $ws.Range("N$r").Value = $item.Price              # N - price
    $ws.Range("P$r").Value = $item.Sale               # P - sale price
    $ws.Range("Q$r").Value = $item.Trans              # Q - number of transactions
}
